# RQMT_game.docx edit
# - Colors several existing TODO bullets red (RGB FF0000 == wdColor 255)
# - Rewrites the "Splash screen" bullet text (adds "/ lightning", swaps the
#   bracketed tile names)
# - Adds four new TODO bullets (ROTATING / WORM / SEND / MATRIX) and expands
#   the old "Big lightning strikes when sending" bullet into a DESTRUCTION
#   bullet that still contains that original sentence, braced.

$d = $word.ActiveDocument
$wdColorRed = 255

# ---------------------------------------------------------------------------
# Helper: locate the paragraph whose text starts with $prefix. Re-run after
# any structural edit (paragraph insert/delete) since indices shift.
# ---------------------------------------------------------------------------
function Find-ParagraphIndex($prefix) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text.StartsWith($prefix)) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1. "Move worm management to GameTable class"
#    -> whole paragraph (every run + the paragraph mark) turns red.
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item((Find-ParagraphIndex("Move worm management to")))
$p.Range.Font.Color = $wdColorRed

# ---------------------------------------------------------------------------
# 2. "Bonus items { CHARGES, bombs, clocks } put them into GameTable class"
#    -> only the leading "Bonus items" text turns red; everything else (and
#    the paragraph mark) is untouched.
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item((Find-ParagraphIndex("Bonus items { CHARGES")))
$s = $p.Range.Start
$lead = "Bonus items"
$d.Range($s, $s + $lead.Length).Font.Color = $wdColorRed

# ---------------------------------------------------------------------------
# 3. "Splash screen with tiles { PIPE CHARGER }"
#    -> "Splash screen with tiles / lightning { ELECTRO PIPES }"
#    -> whole paragraph (text + paragraph mark) turns red.
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item((Find-ParagraphIndex("Splash screen with tiles")))
$p.Range.Find.Execute("PIPE CHARGER", $false, $false, $false, $false, $false, `
    $true, 1, $false, "ELECTRO PIPES", 2)
$p.Range.Find.Execute("tiles { ELECTRO", $false, $false, $false, $false, $false, `
    $true, 1, $false, "tiles / lightning { ELECTRO", 2)
$p.Range.Font.Color = $wdColorRed

# ---------------------------------------------------------------------------
# 4. "Particle generators { with time limit }" -> whole paragraph turns red.
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item((Find-ParagraphIndex("Particle generators")))
$p.Range.Font.Color = $wdColorRed

# ---------------------------------------------------------------------------
# 5. "Lightning updates instead of re-generate" -> whole paragraph turns red.
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item((Find-ParagraphIndex("Lightning updates instead of re-generate")))
$p.Range.Font.Color = $wdColorRed

# ---------------------------------------------------------------------------
# 6. Expand "Big lightning strikes when sending" into five bullets:
#       ROTATING (red), WORM (no color), SEND (red),
#       DESTRUCTION (red, contains the original sentence braced),
#       MATRIX (no color)
# ---------------------------------------------------------------------------

# Three new empty bullets directly before the old one (ROTATING/WORM/SEND).
$p = $d.Paragraphs.Item((Find-ParagraphIndex("Big lightning strikes when sending")))
$p.Range.InsertParagraphBefore()
$p = $d.Paragraphs.Item((Find-ParagraphIndex("Big lightning strikes when sending")))
$p.Range.InsertParagraphBefore()
$p = $d.Paragraphs.Item((Find-ParagraphIndex("Big lightning strikes when sending")))
$p.Range.InsertParagraphBefore()

# One new empty bullet directly after the old one (MATRIX).
$p = $d.Paragraphs.Item((Find-ParagraphIndex("Big lightning strikes when sending")))
$p.Range.InsertParagraphAfter()

$destIdx = Find-ParagraphIndex("Big lightning strikes when sending")
$rotIdx = $destIdx - 3
$wormIdx = $destIdx - 2
$sendIdx = $destIdx - 1
$matrixIdx = $destIdx + 1

$pRot = $d.Paragraphs.Item($rotIdx)
$pRot.Range.Text = "ROTATING -- tap & drag tile to rotate (also actually rotating a tile)"
$pRot.Range.Font.Color = $wdColorRed

$pWorm = $d.Paragraphs.Item($wormIdx)
$pWorm.Range.Text = "WORM -- make it destructible, make it smarter, make it move smoother"

$pSend = $d.Paragraphs.Item($sendIdx)
$pSend.Range.Text = "SEND -- when connecting L-R, you have 48 frames to rotate tiles. if nothing happens, game auto-sends. meanwhile the connecting lightning changes color"
$pSend.Range.Font.Color = $wdColorRed

$pMatrix = $d.Paragraphs.Item($matrixIdx)
$pMatrix.Range.Text = "MATRIX -- [3D] should have depth planes, big letters to small letters -- also you should scroll horizontally through it"

$pDest = $d.Paragraphs.Item($destIdx)
$destStart = $pDest.Range.Start
$d.Range($destStart, $destStart).InsertBefore("DESTRUCTION -- when tiles disappear (bomb, arrow, send), some lightning appears in place for a bit { ")
# Text now reads "...in place for a bit { Big lightning strikes when sending";
# close the brace right before the paragraph mark.
$pDest = $d.Paragraphs.Item($destIdx)
$destEnd = $pDest.Range.End - 1
$d.Range($destEnd, $destEnd).InsertBefore(" }")
$pDest = $d.Paragraphs.Item($destIdx)
$pDest.Range.Font.Color = $wdColorRed

Write-Output "done"
